# Updates cryptos list: apply new Price (D) and Volume(1h) (E) values per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "54.314.19"
Set-TextValue $ws.Range("E2") "  +0.63%  "
Set-TextValue $ws.Range("D3") "2.270.18"
Set-TextValue $ws.Range("E3") "  +0.11%  "
Set-TextValue $ws.Range("E4") "  +0.24%  "
Set-TextValue $ws.Range("D5") "499.33"
Set-TextValue $ws.Range("E5") "  +1.37%  "
Set-TextValue $ws.Range("D6") "129.15"
Set-TextValue $ws.Range("E6") "  +1.59%  "
Set-TextValue $ws.Range("E7") "  +0.04%  "
Set-TextValue $ws.Range("E8") "  -0.13%  "
Set-TextValue $ws.Range("E9") "  +0.35%  "
Set-TextValue $ws.Range("E10") "  +0.82%  "
Set-TextValue $ws.Range("D11") "0.336"
Set-TextValue $ws.Range("E11") "  +3.73%  "
Set-TextValue $ws.Range("D12") "4.91"
Set-TextValue $ws.Range("E12") "  +5.76%  "
Set-TextValue $ws.Range("D13") "23.16"
Set-TextValue $ws.Range("E13") "  +6.77%  "
Set-TextValue $ws.Range("D14") "2.672.86"
Set-TextValue $ws.Range("E14") "  +0.91%  "
Set-TextValue $ws.Range("D15") "54.291.35"
Set-TextValue $ws.Range("E15") "  +0.74%  "
Set-TextValue $ws.Range("E16") "  +0.94%  "
Set-TextValue $ws.Range("D17") "2.280.74"
Set-TextValue $ws.Range("E17") "  +1.25%  "
Set-TextValue $ws.Range("E18") "  +2.83%  "
Set-TextValue $ws.Range("E19") "  +1.79%  "
Set-TextValue $ws.Range("D20") "304.64"
Set-TextValue $ws.Range("E20") "  +1.84%  "
Set-TextValue $ws.Range("E21") "  -1.27%  "
Set-TextValue $ws.Range("E22") "  +0.17%  "
Set-TextValue $ws.Range("D23") "60.52"
Set-TextValue $ws.Range("E23") "  -2.05%  "
Set-TextValue $ws.Range("D24") "0.997"
Set-TextValue $ws.Range("E24") "  -2.05%  "
Set-TextValue $ws.Range("E25") "  +1.08%  "
Set-TextValue $ws.Range("D26") "7.35"
Set-TextValue $ws.Range("E26") "  +4.70%  "
Set-TextValue $ws.Range("D27") "175.27"
Set-TextValue $ws.Range("E27") "  +4.94%  "
Set-TextValue $ws.Range("D28") "0.0₃0702"
Set-TextValue $ws.Range("E28") "  +3.01%  "
Set-TextValue $ws.Range("D29") "6.02"
Set-TextValue $ws.Range("E29") "  +3.18%  "
Set-TextValue $ws.Range("E30") "  +0.57%  "
Set-TextValue $ws.Range("E32") "  +0.02%  "
Set-TextValue $ws.Range("E33") "  +1.39%  "
Set-TextValue $ws.Range("E34") "  -0.10%  "
Set-TextValue $ws.Range("E35") "  +5.95%  "
Set-TextValue $ws.Range("E36") "  +1.97%  "
Set-TextValue $ws.Range("E37") "  +1.74%  "
Set-TextValue $ws.Range("E38") "  +1.27%  "
Set-TextValue $ws.Range("D39") "1.40"
Set-TextValue $ws.Range("E39") "  +0.57%  "
Set-TextValue $ws.Range("D40") "3.38"
Set-TextValue $ws.Range("E40") "  +1.16%  "
Set-TextValue $ws.Range("D41") "4.83"
Set-TextValue $ws.Range("E41") "  +0.74%  "
Set-TextValue $ws.Range("D42") "125.12"
Set-TextValue $ws.Range("E42") "  +0.11%  "
Set-TextValue $ws.Range("D43") "0.0491"
Set-TextValue $ws.Range("E43") "  +2.00%  "
Set-TextValue $ws.Range("E44") "  +1.31%  "
Set-TextValue $ws.Range("D45") "245.54"
Set-TextValue $ws.Range("E45") "  +4.05%  "
Set-TextValue $ws.Range("E46") "  +1.25%  "
Set-TextValue $ws.Range("E47") "  +1.56%  "
Set-TextValue $ws.Range("E48") "  +1.86%  "
Set-TextValue $ws.Range("E49") "  +0.88%  "
Set-TextValue $ws.Range("D50") "16.20"
Set-TextValue $ws.Range("E50") "  +0.99%  "
Set-TextValue $ws.Range("D51") "1.52"
Set-TextValue $ws.Range("E51") "  +3.04%  "
